# fix for geography data structure
#
# 1) Two "general" horizontal-alignment cell styles (the default style used
#    for the data columns on both sheets, including cell C1 on the
#    "Geography Data" sheet) should be left-aligned instead.
# 2) The data rows 2-20 on the "Asset Class Data" sheet get a bit taller
#    (18.75pt instead of 18pt / 17.25pt).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Asset Class Data")
$ws2 = $wb.Worksheets.Item("Geography Data")

$xlLeft = -4131

# --- Row heights: rows 2 through 20 on "Asset Class Data" -> 18.75pt ---
$ws1.Range("A2:B20").EntireRow.RowHeight = 18.75

# --- Alignment: switch "general" horizontal alignment to "left" ---
$ws1.Columns("A:B").HorizontalAlignment = $xlLeft
$ws2.Columns("A:C").HorizontalAlignment = $xlLeft
